$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price column cells whose new values look numeric,
# so Excel does not auto-convert them to numbers (losing formatting like trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.105.51"
$ws.Range("D3").Value = "1.780.30"
$ws.Range("D5").Value = "225.83"
$ws.Range("D8").Value = "31.67"
$ws.Range("D9").Value = "0.292"
$ws.Range("D10").Value = "0.0689"
$ws.Range("D12").Value = "2.035.25"
$ws.Range("D13").Value = "1.793.26"
$ws.Range("D14").Value = "10.90"
$ws.Range("D15").Value = "34.091.46"
$ws.Range("D16").Value = "0.620"
$ws.Range("D18").Value = "67.82"
$ws.Range("D19").Value = "245.42"
$ws.Range("D20").Value = "0.0₃0793"
$ws.Range("D21").Value = "11.02"
$ws.Range("D23").Value = "4.10"
$ws.Range("D25").Value = "161.38"
$ws.Range("D26").Value = "7.17"
$ws.Range("D31").Value = "0.0519"
$ws.Range("D32").Value = "3.72"
$ws.Range("D35").Value = "1.439.34"
$ws.Range("D36").Value = "0.656"
$ws.Range("D40").Value = "80.15"
$ws.Range("D42").Value = "0.924"
$ws.Range("D44").Value = "13.39"
$ws.Range("D45").Value = "6.08"
$ws.Range("D46").Value = "0.0509"
$ws.Range("D48").Value = "0.0₆0137"
$ws.Range("D49").Value = "1.938.55"
$ws.Range("D50").Value = "104.40"
$ws.Range("D51").Value = "1.00"

# Restore default style (keeps text value, drops the explicit text number format)
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

# Remaining (non-numeric-risk) cell updates
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("E21").Value = "  +3.80%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("E33").Value = "  +5.21%  "
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("E37").Value = "  +5.83%  "
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  +1.76%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E45").Value = "  +3.47%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("E51").Value = "  +0.13%  "
